$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Chaolin%Huang%NULL%14,  Yeming%Wang%NULL%2,  Xingwang%Li%NULL%1,  Lili%Ren%NULL%1,  Jianping%Zhao%NULL%1,  Yi%Hu%NULL%3,  Li%Zhang%NULL%3,  Guohui%Fan%NULL%2,  Jiuyang%Xu%NULL%2,  Xiaoying%Gu%NULL%2,  Zhenshun%Cheng%NULL%2,  Ting%Yu%NULL%4,  Jiaan%Xia%NULL%1,  Yuan%Wei%NULL%3,  Wenjuan%Wu%NULL%1,  Xuelei%Xie%NULL%1,  Wen%Yin%NULL%1,  Hui%Li%NULL%3,  Min%Liu%NULL%1,  Yan%Xiao%NULL%2,  Hong%Gao%NULL%2,  Li%Guo%NULL%2,  Jungang%Xie%NULL%1,  Guangfa%Wang%NULL%1,  Rongmeng%Jiang%NULL%1,  Zhancheng%Gao%NULL%1,  Qi%Jin%NULL%2,  Jianwei%Wang%wangjw28@163.com%1,  Bin%Cao%caobin_ben@163.com%3]'
$ws.Range("E3").Value = '[Jasper Fuk-Woo%Chan%NULL%3,  Shuofeng%Yuan%NULL%1,  Kin-Hang%Kok%NULL%1,  Kelvin Kai-Wang%To%NULL%1,  Hin%Chu%NULL%1,  Jin%Yang%NULL%1,  Fanfan%Xing%NULL%1,  Jieling%Liu%NULL%1,  Cyril Chik-Yan%Yip%NULL%1,  Rosana Wing-Shan%Poon%NULL%1,  Hoi-Wah%Tsoi%NULL%1,  Simon Kam-Fai%Lo%NULL%1,  Kwok-Hung%Chan%NULL%1,  Vincent Kwok-Man%Poon%NULL%1,  Wan-Mui%Chan%NULL%1,  Jonathan Daniel%Ip%NULL%1,  Jian-Piao%Cai%NULL%1,  Vincent Chi-Chung%Cheng%NULL%1,  Honglin%Chen%NULL%1,  Christopher Kim-Ming%Hui%NULL%1,  Kwok-Yung%Yuen%NULL%1]'
$ws.Range("E4").Value = '[Qun%Li%NULL%2,  Xuhua%Guan%NULL%1,  Peng%Wu%NULL%1,  Xiaoye%Wang%NULL%1,  Lei%Zhou%NULL%1,  Yeqing%Tong%NULL%1,  Ruiqi%Ren%NULL%1,  Kathy S.M.%Leung%NULL%1,  Eric H.Y.%Lau%NULL%1,  Jessica Y.%Wong%NULL%1,  Xuesen%Xing%NULL%1,  Nijuan%Xiang%NULL%1,  Yang%Wu%NULL%1,  Chao%Li%NULL%1,  Qi%Chen%NULL%1,  Dan%Li%NULL%1,  Tian%Liu%NULL%1,  Jing%Zhao%NULL%1,  Man%Liu%NULL%1,  Wenxiao%Tu%NULL%1,  Chuding%Chen%NULL%1,  Lianmei%Jin%NULL%1,  Rui%Yang%NULL%1,  Qi%Wang%NULL%1,  Suhua%Zhou%NULL%1,  Rui%Wang%NULL%1,  Hui%Liu%NULL%1,  Yinbo%Luo%NULL%1,  Yuan%Liu%NULL%1,  Ge%Shao%NULL%1,  Huan%Li%NULL%1,  Zhongfa%Tao%NULL%1,  Yang%Yang%NULL%1,  Zhiqiang%Deng%NULL%1,  Boxi%Liu%NULL%1,  Zhitao%Ma%NULL%1,  Yanping%Zhang%NULL%1,  Guoqing%Shi%NULL%1,  Tommy T.Y.%Lam%NULL%1,  Joseph T.%Wu%NULL%1,  George F.%Gao%NULL%1,  Benjamin J.%Cowling%NULL%1,  Bo%Yang%NULL%2,  Bo%Yang%NULL%0,  Gabriel M.%Leung%NULL%1,  Zijian%Feng%NULL%1]'
$ws.Range("E5").Value = '[Nanshan%Chen%NULL%9,  Min%Zhou%NULL%1,  Xuan%Dong%NULL%1,  Jieming%Qu%NULL%2,  Fengyun%Gong%NULL%1,  Yang%Han%NULL%1,  Yang%Qiu%NULL%1,  Jingli%Wang%NULL%1,  Ying%Liu%NULL%2,  Yuan%Wei%NULL%0,  Jia''an%Xia%NULL%2,  Ting%Yu%NULL%0,  Xinxin%Zhang%NULL%1,  Li%Zhang%NULL%0]'
$ws.Range("E6").Value = '[Fengxiang%Song%NULL%2,  Nannan%Shi%NULL%2,  Nannan%Shi%NULL%0,  Fei%Shan%NULL%1,  Zhiyong%Zhang%NULL%1,  Jie%Shen%NULL%1,  Hongzhou%Lu%NULL%2,  Yun%Ling%NULL%2,  Yebin%Jiang%NULL%2,  Yebin%Jiang%NULL%0,  Yuxin%Shi%shiyuxin@shphc.org.cn%3,  Yuxin%Shi%shiyuxin@shphc.org.cn%0]'
$ws.Range("E7").Value = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%11,  Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%1,  Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%1,  Subert-Salas%Lizandra%coreGivesNoEmail%1]'
$ws.Range("E8").Value = '[Kui%Liu%NULL%6,  Yuan-Yuan%Fang%NULL%1,  Yan%Deng%NULL%2,  Wei%Liu%NULL%1,  Mei-Fang%Wang%NULL%1,  Jing-Ping%Ma%NULL%1,  Wei%Xiao%NULL%1,  Ying-Nan%Wang%NULL%1,  Min-Hua%Zhong%NULL%1,  Cheng-Hong%Li%NULL%1,  Guang-Cai%Li%NULL%1,  Hui-Guo%Liu%NULL%2,  Xiu-Yuan%Hao%NULL%4,  Xiu-Yuan%Hao%NULL%0,  Pei-Fang%Wei%NULL%4]'
$ws.Range("E9").Value = '[De%Chang%xref no email%2,  Minggui%Lin%xref no email%1,  Lai%Wei%xref no email%1,  Lixin%Xie%xref no email%1,  Guangfa%Zhu%xref no email%1,  Charles S.%Dela Cruz%xref no email%1,  Lokesh%Sharma%xref no email%1]'
$ws.Range("I9").Value = '_CROSSREF'
$ws.Range("E10").Value = '[ Yingxia%Liu%null%2,   Yang%Yang%null%1,   Cong%Zhang%null%1,   Fengming%Huang%null%1,   Fuxiang%Wang%null%1,   Jing%Yuan%null%1,   Zhaoqin%Wang%null%2,   Jinxiu%Li%null%1,   Jianming%Li%null%1,   Cheng%Feng%null%1,   Zheng%Zhang%null%1,   Lifei%Wang%null%1,   Ling%Peng%null%1,   Li%Chen%null%1,   Yuhao%Qin%null%1,   Dandan%Zhao%null%1,   Shuguang%Tan%null%1,   Lu%Yin%null%1,   Jun%Xu%null%1,   Congzhao%Zhou%null%1,   Chengyu%Jiang%null%1,   Lei%Liu%null%2]'
$ws.Range("E11").Value = '[Li-Li%Ren%NULL%3,  Ye-Ming%Wang%NULL%1,  Zhi-Qiang%Wu%NULL%1,  Zi-Chun%Xiang%NULL%1,  Li%Guo%NULL%0,  Teng%Xu%NULL%1,  Yong-Zhong%Jiang%NULL%1,  Yan%Xiong%NULL%1,  Yong-Jun%Li%NULL%1,  Xing-Wang%Li%NULL%1,  Hui%Li%NULL%0,  Guo-Hui%Fan%NULL%1,  Xiao-Ying%Gu%NULL%1,  Yan%Xiao%NULL%0,  Hong%Gao%NULL%0,  Jiu-Yang%Xu%NULL%1,  Fan%Yang%NULL%1,  Xin-Ming%Wang%NULL%1,  Chao%Wu%NULL%1,  Lan%Chen%NULL%1,  Yi-Wei%Liu%NULL%1,  Bo%Liu%NULL%1,  Jian%Yang%NULL%1,  Xiao-Rui%Wang%NULL%1,  Jie%Dong%NULL%1,  Li%Li%NULL%2,  Chao-Lin%Huang%NULL%1,  Jian-Ping%Zhao%NULL%1,  Yi%Hu%NULL%0,  Zhen-Shun%Cheng%NULL%1,  Lin-Lin%Liu%NULL%1,  Zhao-Hui%Qian%NULL%1,  Chuan%Qin%NULL%1,  Qi%Jin%NULL%0,  Bin%Cao%NULL%0,  Jian-Wei%Wang%NULL%1,  Xiu-Yuan%Hao%NULL%0,  Xiu-Yuan%Hao%NULL%0,  Pei-Fang%Wei%NULL%0]'
$ws.Range("E12").Value = '[Ning%Tang%NULL%4,  Dengju%Li%NULL%1,  Xiong%Wang%NULL%1,  Ziyong%Sun%NULL%3]'
$ws.Range("E13").Value = '[Akdis%Cezmi A%coreGivesNoEmail%9,  Cao%Yi-yuan%coreGivesNoEmail%1,  Dong%Xiang%coreGivesNoEmail%1,  Gao%Ya-dong%coreGivesNoEmail%1,  Yan%You-qin%coreGivesNoEmail%1,  Yang%Yi-bin%coreGivesNoEmail%1,  Yuan%Ya-dong%coreGivesNoEmail%1,  Zhang%Jin-jin%coreGivesNoEmail%1]'
$ws.Range("E14").Value = '[Xiao-Wei%Xu%NULL%7,  Xiao-Xin%Wu%NULL%1,  Xian-Gao%Jiang%NULL%1,  Kai-Jin%Xu%NULL%1,  Ling-Jun%Ying%NULL%1,  Chun-Lian%Ma%NULL%1,  Shi-Bo%Li%NULL%1,  Hua-Ying%Wang%NULL%1,  Sheng%Zhang%NULL%1,  Hai-Nv%Gao%NULL%1,  Ji-Fang%Sheng%NULL%1,  Hong-Liu%Cai%NULL%1,  Yun-Qing%Qiu%NULL%1,  Lan-Juan%Li%NULL%1]'
$ws.Range("E15").Value = '[Jiong%Wu%NULL%1,  Xiaojia%Wu%NULL%1,  Wenbing%Zeng%NULL%1,  Dajing%Guo%NULL%2,  Zheng%Fang%NULL%2,  Linli%Chen%NULL%2,  Huizhe%Huang%NULL%1,  Chuanming%Li%NULL%2]'
$ws.Range("E16").Value = '[Kin On%Kwok%kkokwok@cuhk.edu.hk%2,  Valerie%Wong%NULL%1,  Vivian Wan In%Wei%NULL%1,  Samuel Yeung Shan%Wong%NULL%1,  Julian Wei-Tze%Tang%julian.tang@uhl-tr.nhs.uk%1]'
$ws.Range("E17").Value = '[Heshui%Shi%NULL%2,  Xiaoyu%Han%NULL%1,  Nanchuan%Jiang%NULL%1,  Yukun%Cao%NULL%1,  Osamah%Alwalid%NULL%1,  Jin%Gu%NULL%1,  Yanqing%Fan%1024932023@qq.com%2,  Chuansheng%Zheng%hqzcsxh@sina.com%1]'
$ws.Range("E18").Value = '[Xiaobo%Yang%NULL%16,  Yuan%Yu%NULL%1,  Jiqian%Xu%NULL%1,  Huaqing%Shu%NULL%1,  Jia''an%Xia%NULL%0,  Hong%Liu%NULL%1,  Yongran%Wu%NULL%1,  Lu%Zhang%NULL%1,  Zhui%Yu%NULL%1,  Minghao%Fang%NULL%1,  Ting%Yu%NULL%0,  Yaxin%Wang%NULL%1,  Shangwen%Pan%NULL%1,  Xiaojing%Zou%NULL%1,  Shiying%Yuan%NULL%1,  You%Shang%NULL%1]'
$ws.Range("E19").Value = '[Yu-Huan%Xu%NULL%2,  Jing-Hui%Dong%NULL%1,  Wei-Min%An%NULL%1,  Xiao-Yan%Lv%15001008285@139.com%1,  Xiao-Ping%Yin%yinxiaoping78@sina.com%1,  Jian-Zeng%Zhang%NULL%1,  Li%Dong%NULL%1,  Xi%Ma%NULL%1,  Hong-Jie%Zhang%NULL%1,  Bu-Lang%Gao%browngao@163.com%1]'
$ws.Range("E20").Value = '[Tao%Ai%NULL%3,  Zhenlu%Yang%NULL%2,  Zhenlu%Yang%NULL%0,  Hongyan%Hou%NULL%2,  Hongyan%Hou%NULL%0,  Chenao%Zhan%NULL%1,  Chong%Chen%NULL%1,  Wenzhi%Lv%NULL%2,  Wenzhi%Lv%NULL%0,  Qian%Tao%NULL%2,  Qian%Tao%NULL%0,  Ziyong%Sun%NULL%0,  Ziyong%Sun%NULL%0,  Liming%Xia%xialiming2017@outlook.com%3]'
$ws.Range("E21").Value = '[Wenjie%Yang%NULL%4,  Qiqi%Cao%NULL%2,  Le%Qin%NULL%2,  Xiaoyang%Wang%NULL%1,  Zenghui%Cheng%NULL%1,  Ashan%Pan%NULL%1,  Jianyi%Dai%NULL%1,  Qingfeng%Sun%NULL%1,  Fengquan%Zhao%NULL%1,  Jieming%Qu%NULL%0,  Fuhua%Yan%NULL%2]'
$ws.Range("E22").Value = '[Sijia%Tian%NULL%2,  Nan%Hu%NULL%1,  Jing%Lou%NULL%1,  Kun%Chen%NULL%1,  Xuqin%Kang%NULL%1,  Zhenjun%Xiang%NULL%1,  Hui%Chen%NULL%1,  Dali%Wang%NULL%1,  Ning%Liu%NULL%1,  Dong%Liu%NULL%1,  Gang%Chen%NULL%1,  Yongliang%Zhang%NULL%1,  Dou%Li%NULL%1,  Jianren%Li%NULL%1,  Huixin%Lian%NULL%1,  Shengmei%Niu%NULL%1,  Luxi%Zhang%NULL%1,  Jinjun%Zhang%NULL%1]'
$ws.Range("E23").Value = '[Yihui%Huang%NULL%3,  Mengqi%Tu%NULL%1,  Shipei%Wang%NULL%1,  Sichao%Chen%NULL%1,  Wei%Zhou%NULL%1,  Danyang%Chen%NULL%1,  Lin%Zhou%NULL%1,  Min%Wang%NULL%1,  Yan%Zhao%NULL%1,  Wen%Zeng%NULL%1,  Qi%Huang%NULL%1,  Hai''bo%Xu%NULL%1,  Zeming%Liu%NULL%1,  Liang%Guo%NULL%1]'
$ws.Range("E24").Value = '[Xi%Xu%NULL%4,  Chengcheng%Yu%NULL%1,  Jing%Qu%NULL%1,  Lieguang%Zhang%NULL%1,  Songfeng%Jiang%NULL%1,  Deyang%Huang%NULL%1,  Bihua%Chen%NULL%1,  Zhiping%Zhang%NULL%1,  Wanhua%Guan%NULL%1,  Zhoukun%Ling%NULL%1,  Rui%Jiang%NULL%1,  Tianli%Hu%NULL%1,  Yan%Ding%NULL%1,  Lin%Lin%NULL%1,  Qingxin%Gan%NULL%1,  Liangping%Luo%tluolp@jnu.edu.cn%1,  Xiaoping%Tang%xtang@21cn.com%1,  Jinxin%Liu%Liujx83710378@126.com%1]'
$ws.Range("E25").Value = '[Wei%Liu%NULL%8,  Zhao-Wu%Tao%NULL%1,  Lei%Wang%NULL%1,  Ming-Li%Yuan%NULL%1,  Kui%Liu%NULL%1,  Ling%Zhou%NULL%1,  Shuang%Wei%NULL%1,  Yan%Deng%NULL%0,  Jing%Liu%NULL%1,  Hui-Guo%Liu%NULL%0,  Ming%Yang%NULL%1,  Yi%Hu%NULL%0,  Pei-Fang%Wei%NULL%0,  Pei-Fang%Wei%NULL%0]'
$ws.Range("E26").Value = '[Kunhua%Li%NULL%3,  Jiong%Wu%NULL%1,  Faqi%Wu%NULL%1,  Dajing%Guo%NULL%0,  Linli%Chen%NULL%0,  Zheng%Fang%NULL%0,  Chuanming%Li%NULL%0]'
$ws.Range("E27").Value = '[Jian%Wu%NULL%2,  Jun%Liu%NULL%2,  Xinguo%Zhao%NULL%1,  Chengyuan%Liu%NULL%1,  Wei%Wang%NULL%1,  Dawei%Wang%NULL%1,  Wei%Xu%NULL%1,  Chunyu%Zhang%NULL%1,  Jiong%Yu%NULL%1,  Bin%Jiang%NULL%1,  Hongcui%Cao%hccao@zju.edu.cn%1,  Lanjuan%Li%NULL%1]'
$ws.Range("E28").Value = '[Jianlei%Cao%NULL%12,  Xiaorong%Hu%NULL%1,  Wenlin%Cheng%NULL%1,  Lei%Yu%NULL%1,  Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%1,  Qiang%Liu%liuqiang@irm-cams.ac.cn%2,  Qiang%Liu%liuqiang@irm-cams.ac.cn%0]'
$ws.Range("E29").Value = '[Qiurong%Ruan%NULL%5,  Kun%Yang%NULL%2,  Kun%Yang%NULL%0,  Wenxia%Wang%NULL%2,  Wenxia%Wang%NULL%0,  Lingyu%Jiang%NULL%2,  Lingyu%Jiang%NULL%0,  Jianxin%Song%songsingsjx@sina.com%2,  Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range("E30").Value = '[Wei%Zhao%xref no email%1,  Zheng%Zhong%xref no email%1,  Xingzhi%Xie%xref no email%1,  Qizhi%Yu%xref no email%1,  Jun%Liu%xref no email%0]'
$ws.Range("I30").Value = '_CROSSREF'
$ws.Range("E31").Value = '[Fei%Xiao%NULL%1,  Meiwen%Tang%NULL%1,  Xiaobin%Zheng%NULL%1,  Ye%Liu%ly77219@163.com%1,  Xiaofeng%Li%zdwylxf@163.com%1,  Hong%Shan%shanhong@mail.sysu.edu.cn%1]'
$ws.Range("E32").Value = '[Ying%Xiong%NULL%1,  Dong%Sun%NULL%1,  Yao%Liu%NULL%1,  Yanqing%Fan%NULL%0,  Lingyun%Zhao%NULL%1,  Xiaoming%Li%NULL%1,  Wenzhen%Zhu%NULL%1]'
$ws.Range("E33").Value = '[Barnaby Edward%Young%xref no email%6,  Sean Wei Xiang%Ong%xref no email%1,  Shirin%Kalimuddin%xref no email%1,  Jenny G.%Low%xref no email%1,  Seow Yen%Tan%xref no email%1,  Jiashen%Loh%xref no email%1,  Oon-Tek%Ng%xref no email%1,  Kalisvar%Marimuthu%xref no email%1,  Li Wei%Ang%xref no email%1,  Tze Minn%Mak%xref no email%1,  Sok Kiang%Lau%xref no email%1,  Danielle E.%Anderson%xref no email%1,  Kian Sing%Chan%xref no email%1,  Thean Yen%Tan%xref no email%1,  Tong Yong%Ng%xref no email%1,  Lin%Cui%xref no email%1,  Zubaidah%Said%xref no email%1,  Lalitha%Kurupatham%xref no email%1,  Mark I-Cheng%Chen%xref no email%1,  Monica%Chan%xref no email%1,  Shawn%Vasoo%xref no email%1,  Lin-Fa%Wang%xref no email%1,  Boon Huan%Tan%xref no email%1,  Raymond Tzer Pin%Lin%xref no email%1,  Vernon Jian Ming%Lee%xref no email%1,  Yee-Sin%Leo%xref no email%1,  David Chien%Lye%xref no email%1]'
$ws.Range("I33").Value = '_CROSSREF'
$ws.Range("E34").Value = '[JingCheng%Zhang%zjc1983@126.com%1,  SaiBin%Wang%NULL%2,  SaiBin%Wang%NULL%0,  YaDong%Xue%NULL%1]'
$ws.Range("E35").Value = '[Yan%Li%xref no email%1,  Liming%Xia%xref no email%0]'
$ws.Range("I35").Value = '_CROSSREF'
$ws.Range("E36").Value = '[Zhiliang%Hu%NULL%2,  Ci%Song%NULL%1,  Chuanjun%Xu%NULL%1,  Guangfu%Jin%NULL%1,  Yaling%Chen%NULL%1,  Xin%Xu%NULL%1,  Hongxia%Ma%NULL%1,  Wei%Chen%NULL%1,  Yuan%Lin%NULL%1,  Yishan%Zheng%NULL%1,  Jianming%Wang%NULL%1,  Zhibin%Hu%zhibin_hu@njmu.edu.cn%1,  Yongxiang%Yi%ian0126@126.com%1,  Hongbing%Shen%NULL%1]'
$ws.Range("E37").Value = '[Yang%Yao%NULL%1,  Yao%Tian%NULL%2,  Yao%Tian%NULL%0,  Jing%Zhou%NULL%1,  Xuan%Ma%NULL%1,  Min%Yang%NULL%1,  ShengYu%Wang%NULL%1]'
$ws.Range("E38").Value = '[Shuchang%Zhou%xref no email%1,  Yujin%Wang%xref no email%1,  Tingting%Zhu%xref no email%1,  Liming%Xia%xref no email%0]'
$ws.Range("I38").Value = '_CROSSREF'
$ws.Range("E39").Value = '[Gianfranco%Spiteri%NULL%1,  James%Fielding%NULL%1,  Michaela%Diercke%NULL%1,  Christine%Campese%NULL%1,  Vincent%Enouf%NULL%1,  Alexandre%Gaymard%NULL%1,  Antonino%Bella%NULL%1,  Paola%Sognamiglio%NULL%1,  Maria José%Sierra Moros%NULL%1,  Antonio Nicolau%Riutort%NULL%1,  Yulia V.%Demina%NULL%1,  Romain%Mahieu%NULL%1,  Markku%Broas%NULL%1,  Malin%Bengnér%NULL%1,  Silke%Buda%NULL%1,  Julia%Schilling%NULL%1,  Laurent%Filleul%NULL%1,  Agnès%Lepoutre%NULL%1,  Christine%Saura%NULL%1,  Alexandra%Mailles%NULL%1,  Daniel%Levy-Bruhl%NULL%1,  Bruno%Coignard%NULL%1,  Sibylle%Bernard-Stoecklin%NULL%1,  Sylvie%Behillil%NULL%1,  Sylvie%van der Werf%NULL%1,  Martine%Valette%NULL%1,  Bruno%Lina%NULL%1,  Flavia%Riccardo%NULL%1,  Emanuele%Nicastri%NULL%1,  Inmaculada%Casas%NULL%1,  Amparo%Larrauri%NULL%1,  Magdalena%Salom Castell%NULL%1,  Francisco%Pozo%NULL%1,  Rinat A.%Maksyutov%NULL%1,  Charlotte%Martin%NULL%1,  Marc%Van Ranst%NULL%1,  Nathalie%Bossuyt%NULL%1,  Lotta%Siira%NULL%1,  Jussi%Sane%NULL%1,  Karin%Tegmark-Wisell%NULL%1,  Maria%Palmérus%NULL%1,  Eeva K.%Broberg%NULL%1,  Julien%Beauté%NULL%1,  Pernille%Jorgensen%NULL%1,  Nick%Bundle%NULL%1,  Dmitriy%Pereyaslov%NULL%1,  Cornelia%Adlhoch%NULL%1,  Jukka%Pukkila%NULL%1,  Richard%Pebody%NULL%1,  Sonja%Olsen%NULL%1,  Bruno Christian%Ciancio%NULL%1]'
$ws.Range("E40").Value = '[Lei%Wang%NULL%1,  Yong-hua%Gao%NULL%1,  Li-Li%lou%NULL%1,  Guo-Jun%Zhang%NULL%1]'
$ws.Range("E41").Value = '[Simin%Zhang%NULL%1,  Huaqiao%Li%NULL%1,  Songtao%Huang%NULL%1,  Wei%You%NULL%1,  Huaiqiang%Sun%NULL%1]'
$ws.Range("E42").Value = '[Fei%Zhou%NULL%18,  Ting%Yu%NULL%0,  Ronghui%Du%NULL%1,  Guohui%Fan%NULL%0,  Ying%Liu%NULL%0,  Zhibo%Liu%NULL%1,  Jie%Xiang%NULL%1,  Yeming%Wang%NULL%0,  Bin%Song%NULL%1,  Xiaoying%Gu%NULL%0,  Lulu%Guan%NULL%1,  Yuan%Wei%NULL%0,  Hui%Li%NULL%0,  Xudong%Wu%NULL%1,  Jiuyang%Xu%NULL%0,  Shengjin%Tu%NULL%1,  Yi%Zhang%NULL%1,  Hua%Chen%NULL%1,  Bin%Cao%NULL%0]'
$ws.Range("E43").Value = '[Tracy%Tsang%xref no email%1]'
$ws.Range("I43").Value = '_CROSSREF'
$ws.Range("E44").Value = '[Kai%Liu%hmliukai@126.com%2,  Ying%Chen%NULL%1,  Ruzheng%Lin%NULL%1,  Kunyuan%Han%NULL%1]'
$ws.Range("E45").Value = '[Jun%Chen%NULL%4,  Tangkai%Qi%NULL%1,  Li%Liu%NULL%1,  Yun%Ling%NULL%0,  Zhiping%Qian%NULL%1,  Tao%Li%NULL%1,  Feng%Li%NULL%1,  Qingnian%Xu%NULL%1,  Yuyi%Zhang%NULL%1,  Shuibao%Xu%NULL%1,  Zhigang%Song%NULL%1,  Yigang%Zeng%NULL%1,  Yinzhong%Shen%NULL%1,  Yuxin%Shi%NULL%0,  Tongyu%Zhu%NULL%1,  Hongzhou%Lu%NULL%0]'
$ws.Range("E46").Value = '[Dahai%Zhao%NULL%2,  Feifei%Yao%NULL%1,  Lijie%Wang%NULL%1,  Ling%Zheng%NULL%1,  Yongjun%Gao%NULL%1,  Jun%Ye%NULL%1,  Feng%Guo%NULL%1,  Hui%Zhao%NULL%1,  Rongbao%Gao%rongbaogao@163.com%1]'
$ws.Range("E47").Value = '[Qing%Chen%NULL%1,  Bin%Quan%NULL%1,  Xiaoning%Li%NULL%1,  Guangjian%Gao%NULL%1,  Wenqiang%Zheng%NULL%1,  Jun%Zhang%NULL%1,  Zhiyun%Zhang%NULL%1,  Chunsheng%Liu%NULL%1,  Li%Li%NULL%0,  Chenglin%Wang%NULL%1,  Guihua%Zhang%NULL%1,  Jiajia%Li%NULL%1,  Yunhai%Dai%NULL%1,  Jianghua%Yang%yjhpath@163.com%1,  Wenzheng%Han%11418166@zju.edu.cn%1]'
$ws.Range("E48").Value = '[Chaomin%Wu%xref no email%10,  Xiaoyan%Chen%xref no email%1,  Yanping%Cai%xref no email%1,  Jia\u2019an%Xia%xref no email%1,  Xing%Zhou%xref no email%1,  Sha%Xu%xref no email%1,  Hanping%Huang%xref no email%1,  Li%Zhang%xref no email%0,  Xia%Zhou%xref no email%1,  Chunling%Du%xref no email%1,  Yuye%Zhang%xref no email%1,  Juan%Song%xref no email%1,  Sijiao%Wang%xref no email%1,  Yencheng%Chao%xref no email%1,  Zeyong%Yang%xref no email%1,  Jie%Xu%xref no email%1,  Xin%Zhou%xref no email%1,  Dechang%Chen%xref no email%1,  Weining%Xiong%xref no email%1,  Lei%Xu%xref no email%2,  Feng%Zhou%xref no email%1,  Jinjun%Jiang%xref no email%1,  Chunxue%Bai%xref no email%1,  Junhua%Zheng%xref no email%1,  Yuanlin%Song%xref no email%1]'
$ws.Range("I48").Value = '_CROSSREF'
$ws.Range("E49").Value = '[Wanbo%Zhu%NULL%4,  Kai%Xie%NULL%2,  Kai%Xie%NULL%0,  Hui%Lu%NULL%1,  Lei%Xu%bayinhexl@126.com%0,  Shusheng%Zhou%zhouss108@163.com%1,  Shiyuan%Fang%fangshiyuan2008@126.com%1]'
$ws.Range("E50").Value = '[Zenghui%Cheng%xref no email%3,  Yong%Lu%xref no email%1,  Qiqi%Cao%xref no email%0,  Le%Qin%xref no email%0,  Zilai%Pan%xref no email%1,  Fuhua%Yan%xref no email%0,  Wenjie%Yang%xref no email%1]'
$ws.Range("I50").Value = '_CROSSREF'
$ws.Range("E51").Value = '[Pingzheng%Mo%NULL%5,  Yuanyuan%Xing%NULL%1,  Yu%Xiao%NULL%1,  Liping%Deng%NULL%1,  Qiu%Zhao%NULL%1,  Hongling%Wang%NULL%1,  Yong%Xiong%NULL%1,  Zhenshun%Cheng%NULL%0,  Shicheng%Gao%NULL%1,  Ke%Liang%NULL%1,  Mingqi%Luo%NULL%1,  Tielong%Chen%NULL%1,  Shihui%Song%NULL%1,  Zhiyong%Ma%NULL%1,  Xiaoping%Chen%NULL%1,  Ruiying%Zheng%NULL%1,  Qian%Cao%NULL%1,  Fan%Wang%fanndywang@foxmail.com%1,  Yongxi%Zhang%znact1936@126.com%1]'
$ws.Range("E52").Value = '[Zhongliang%Wang%NULL%6,  Bohan%Yang%NULL%1,  Qianwen%Li%NULL%1,  Lu%Wen%NULL%1,  Ruiguang%Zhang%zrg27@163.com%1]'
$ws.Range("E53").Value = '[Yanrong%Wang%123rong@sohu.com%1,  Yingxia%Liu%NULL%0,  Lei%Liu%NULL%1,  Xianfeng%Wang%NULL%1,  Nijuan%Luo%NULL%1,  Ling%Li%NULL%1]'
$ws.Range("E54").Value = '[Rui%Han%xref no email%1,  Lu%Huang%xref no email%1,  Hong%Jiang%xref no email%1,  Jin%Dong%xref no email%1,  Hongfen%Peng%xref no email%1,  Dongyou%Zhang%xref no email%1]'
$ws.Range("I54").Value = '_CROSSREF'
$ws.Range("E55").Value = '[G -Q%Qian%NULL%2,  N -B%Yang%NULL%1,  F%Ding%NULL%1,  A H Y%Ma%NULL%1,  Z -Y%Wang%NULL%2,  Y -F%Shen%NULL%1,  C -W%Shi%NULL%1,  X%Lian%NULL%1,  J -G%Chu%NULL%1,  L%Chen%chxmin@hotmail.com%1,  Z -Y%Wang%NULL%0,  D -W%Ren%NULL%1,  G -X%Li%NULL%1,  X -Q%Chen%chxmin@hotmail.com%1,  H -J%Shen%NULL%1,  X -M%Chen%chxmin@hotmail.com%1]'
$ws.Range("E56").Value = '[Yong%Gao%NULL%3,  Tuantuan%Li%NULL%1,  Mingfeng%Han%NULL%1,  Xiuyong%Li%NULL%1,  Dong%Wu%NULL%1,  Yuanhong%Xu%NULL%1,  Yulin%Zhu%NULL%1,  Yan%Liu%NULL%1,  Xiaowu%Wang%wangxiaowu19880218@126.com%1,  Linding%Wang%wanglinding@ahmu.edu.cn%1]'
